$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 5875.25
$ws.Cells.Item(43, 9).Value = 6000.5
$ws.Cells.Item(43, 11).Value = 6000.5
$ws.Cells.Item(43, 13).Value = -5931.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 2989.149
$ws.Cells.Item(129, 9).Value = 20420.4
$ws.Cells.Item(129, 10).Value = 914
$ws.Cells.Item(129, 11).Value = 61261.2
$ws.Cells.Item(129, 12).Value = 2742
$ws.Cells.Item(129, 13).Value = -56261.2
$ws.Cells.Item(129, 14).Value = -12742

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1469.8649
$ws.Cells.Item(137, 9).Value = 1210.4445
$ws.Cells.Item(137, 11).Value = 3631.3335
$ws.Cells.Item(137, 13).Value = -1081.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3144.2856
$ws.Cells.Item(138, 9).Value = 2081.3809
$ws.Cells.Item(138, 10).Value = 3434.169
$ws.Cells.Item(138, 11).Value = 6244.1427
$ws.Cells.Item(138, 12).Value = 10302.507
$ws.Cells.Item(138, 13).Value = -1104.1427
$ws.Cells.Item(138, 14).Value = -20582.507

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 26889.588
$ws.Cells.Item(32, 9).Value = 9190.306
$ws.Cells.Item(32, 10).Value = 124916.38
$ws.Cells.Item(32, 11).Value = 9190.306
$ws.Cells.Item(32, 12).Value = 124916.38
$ws.Cells.Item(32, 13).Value = -8903.306
$ws.Cells.Item(32, 14).Value = -125490.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 5988
$ws.Cells.Item(44, 10).Value = 5988
$ws.Cells.Item(44, 12).Value = 5988
$ws.Cells.Item(44, 14).Value = -6964

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1164.3235
$ws.Cells.Item(74, 9).Value = 1180.7391
$ws.Cells.Item(74, 10).Value = 1130
$ws.Cells.Item(74, 11).Value = 1180.7391
$ws.Cells.Item(74, 12).Value = 1130
$ws.Cells.Item(74, 13).Value = -306.7391
$ws.Cells.Item(74, 14).Value = -2878

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1164.3235
$ws.Cells.Item(77, 9).Value = 1180.7391
$ws.Cells.Item(77, 10).Value = 1130
$ws.Cells.Item(77, 11).Value = 5903.6955
$ws.Cells.Item(77, 12).Value = 5650
$ws.Cells.Item(77, 13).Value = -1535.6955
$ws.Cells.Item(77, 14).Value = -14386

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 24615.2
$ws.Cells.Item(80, 9).Value = 20000
$ws.Cells.Item(80, 10).Value = 25128
$ws.Cells.Item(80, 11).Value = 20000
$ws.Cells.Item(80, 12).Value = 25128
$ws.Cells.Item(80, 13).Value = -19002
$ws.Cells.Item(80, 14).Value = -27124

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value = 24615.2
$ws.Cells.Item(83, 9).Value = 20000
$ws.Cells.Item(83, 10).Value = 25128
$ws.Cells.Item(83, 11).Value = 60000
$ws.Cells.Item(83, 12).Value = 75384
$ws.Cells.Item(83, 13).Value = -55008
$ws.Cells.Item(83, 14).Value = -85368

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2286.111
$ws.Cells.Item(122, 9).Value = 1867.5
$ws.Cells.Item(122, 11).Value = 5602.5
$ws.Cells.Item(122, 13).Value = -3152.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 201672.1
$ws.Cells.Item(105, 9).Value = 251885
$ws.Cells.Item(105, 11).Value = 251885
$ws.Cells.Item(105, 13).Value = -250138

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2983.347
$ws.Cells.Item(134, 9).Value = 2926.2974
$ws.Cells.Item(134, 10).Value = 3159.25
$ws.Cells.Item(134, 11).Value = 8778.8922
$ws.Cells.Item(134, 12).Value = 9477.75
$ws.Cells.Item(134, 13).Value = -6243.8922
$ws.Cells.Item(134, 14).Value = -14547.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 3043.4546
$ws.Cells.Item(99, 9).Value = 3239.5
$ws.Cells.Item(99, 10).Value = 2999.889
$ws.Cells.Item(99, 11).Value = 3239.5
$ws.Cells.Item(99, 12).Value = 2999.889
$ws.Cells.Item(99, 13).Value = -1741.5
$ws.Cells.Item(99, 14).Value = -5995.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 3043.4546
$ws.Cells.Item(126, 9).Value = 3239.5
$ws.Cells.Item(126, 10).Value = 2999.889
$ws.Cells.Item(126, 11).Value = 9718.5
$ws.Cells.Item(126, 12).Value = 8999.667000000001
$ws.Cells.Item(126, 13).Value = -7248.5
$ws.Cells.Item(126, 14).Value = -13939.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2065.6667
$ws.Cells.Item(134, 9).Value = 1033.3334
$ws.Cells.Item(134, 10).Value = 2581.8333
$ws.Cells.Item(134, 11).Value = 3100.0002
$ws.Cells.Item(134, 12).Value = 7745.499899999999
$ws.Cells.Item(134, 13).Value = -565.0001999999999
$ws.Cells.Item(134, 14).Value = -12815.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1084.6615
$ws.Cells.Item(5, 9).Value = 668.89655
$ws.Cells.Item(5, 10).Value = 1419.5834
$ws.Cells.Item(5, 11).Value = 2006.68965
$ws.Cells.Item(5, 12).Value = 4258.7502
$ws.Cells.Item(5, 13).Value = -1894.68965
$ws.Cells.Item(5, 14).Value = -4482.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1023.3333
$ws.Cells.Item(34, 9).Value = 60
$ws.Cells.Item(34, 10).Value = 1298.5714
$ws.Cells.Item(34, 11).Value = 180
$ws.Cells.Item(34, 12).Value = 3895.7142
$ws.Cells.Item(34, 13).Value = -96
$ws.Cells.Item(34, 14).Value = -4063.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(58, 8).Value = 1435.3334
$ws.Cells.Item(58, 9).Value = 1650
$ws.Cells.Item(58, 10).Value = 1006
$ws.Cells.Item(58, 11).Value = 4950
$ws.Cells.Item(58, 12).Value = 3018
$ws.Cells.Item(58, 13).Value = -4822
$ws.Cells.Item(58, 14).Value = -3274

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 822.58826
$ws.Cells.Item(113, 9).Value = 1339
$ws.Cells.Item(113, 10).Value = 540.9091
$ws.Cells.Item(113, 11).Value = 4017
$ws.Cells.Item(113, 12).Value = 1622.7273
$ws.Cells.Item(113, 13).Value = -1847
$ws.Cells.Item(113, 14).Value = -5962.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1211.78
$ws.Cells.Item(131, 9).Value = 443.85715
$ws.Cells.Item(131, 10).Value = 1415.9114
$ws.Cells.Item(131, 11).Value = 1331.57145
$ws.Cells.Item(131, 12).Value = 4247.7342
$ws.Cells.Item(131, 13).Value = 3708.42855
$ws.Cells.Item(131, 14).Value = -14327.7342

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 1084.6615
$ws.Cells.Item(135, 9).Value = 668.89655
$ws.Cells.Item(135, 10).Value = 1419.5834
$ws.Cells.Item(135, 11).Value = 6020.068950000001
$ws.Cells.Item(135, 12).Value = 12776.2506
$ws.Cells.Item(135, 13).Value = -3485.068950000001
$ws.Cells.Item(135, 14).Value = -17846.2506

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3636.1365
$ws.Cells.Item(132, 9).Value = 2399.8667
$ws.Cells.Item(132, 10).Value = 6285.2856
$ws.Cells.Item(132, 11).Value = 7199.6001
$ws.Cells.Item(132, 12).Value = 18855.8568
$ws.Cells.Item(132, 13).Value = -4669.6001
$ws.Cells.Item(132, 14).Value = -23915.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1756.4706
$ws.Cells.Item(82, 9).Value = 1276.7
$ws.Cells.Item(82, 10).Value = 2441.8572
$ws.Cells.Item(82, 11).Value = 1276.7
$ws.Cells.Item(82, 12).Value = 2441.8572
$ws.Cells.Item(82, 13).Value = -915.7
$ws.Cells.Item(82, 14).Value = -3163.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1756.4706
$ws.Cells.Item(85, 9).Value = 1276.7
$ws.Cells.Item(85, 10).Value = 2441.8572
$ws.Cells.Item(85, 11).Value = 1276.7
$ws.Cells.Item(85, 12).Value = 2441.8572
$ws.Cells.Item(85, 13).Value = -28.70000000000005
$ws.Cells.Item(85, 14).Value = -4937.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3670.5833
$ws.Cells.Item(93, 9).Value = 3784.7778
$ws.Cells.Item(93, 10).Value = 3328
$ws.Cells.Item(93, 11).Value = 3784.7778
$ws.Cells.Item(93, 12).Value = 3328
$ws.Cells.Item(93, 13).Value = -2536.7778
$ws.Cells.Item(93, 14).Value = -5824

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(56, 8).Value = 29264
$ws.Cells.Item(56, 9).Value = 9195
$ws.Cells.Item(56, 11).Value = 9195
$ws.Cells.Item(56, 13).Value = -8481

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 25454.957
$ws.Cells.Item(136, 9).Value = 84258.586
$ws.Cells.Item(136, 10).Value = 5293.7144
$ws.Cells.Item(136, 11).Value = 252775.758
$ws.Cells.Item(136, 12).Value = 15881.1432
$ws.Cells.Item(136, 13).Value = -250225.758
$ws.Cells.Item(136, 14).Value = -20981.1432
